$wb = $excel.ActiveWorkbook

# Hyperlink-style font color used elsewhere in these sheets (RGB 6495ED -> OLE BGR int)
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# zh-cn sheet (row 8: the 8fce654f-... entry now has a handback file/date/error)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Column P (Error Detail) is now wide enough to show the message
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsZh.Range("I8").Value = "8fce654f-3743-4963-9113-71617b5d245c.md"
$wsZh.Range("I8").Font.Underline = $true
$wsZh.Range("I8").Font.Color = $hyperlinkColor

$wsZh.Range("J8").Value = "8fce654f-3743-4963-9113-71617b5d245c.676def500802cd58587700520369b52e9d505bf4.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-29 14:49:13"
$wsZh.Range("P8").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f011fee5b5ab56006ce0190150791d08b8c034a/e2e/8fce654f-3743-4963-9113-71617b5d245c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d54a556a515dd91ff081d6cdb3df3803c073544e/e2e/8fce654f-3743-4963-9113-71617b5d245c.md."

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9f011fee5b5ab56006ce0190150791d08b8c034a/e2e/8fce654f-3743-4963-9113-71617b5d245c.md", "", "", "8fce654f-3743-4963-9113-71617b5d245c.md")

# ---------------------------------------------------------------------------
# de-de sheet (same row 8 gets its own handback file/date/error)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsDe.Range("I8").Value = "8fce654f-3743-4963-9113-71617b5d245c.md"
$wsDe.Range("I8").Font.Underline = $true
$wsDe.Range("I8").Font.Color = $hyperlinkColor

$wsDe.Range("J8").Value = "8fce654f-3743-4963-9113-71617b5d245c.676def500802cd58587700520369b52e9d505bf4.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-29 14:49:21"
$wsDe.Range("P8").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f011fee5b5ab56006ce0190150791d08b8c034a/e2e/8fce654f-3743-4963-9113-71617b5d245c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d54a556a515dd91ff081d6cdb3df3803c073544e/e2e/8fce654f-3743-4963-9113-71617b5d245c.md."

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9f011fee5b5ab56006ce0190150791d08b8c034a/e2e/8fce654f-3743-4963-9113-71617b5d245c.md", "", "", "8fce654f-3743-4963-9113-71617b5d245c.md")
